$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. A new bound is appended below the existing table: a fixed 2050
#    capacity bound (CAP_BND) of 50 GW on the onshore-wind process.
#    Build it first (off row 5's un-shaded style), before the other
#    string edits below, so CAP_BND lands first in the shared-string
#    table.
# ---------------------------------------------------------------------

$ws.Range("B5:F5").Copy()
$ws.Range("B13:F13").PasteSpecial(-4122)

$ws.Range("B13").Value = "UP"
$ws.Range("C13").Value = "CAP_BND"
$ws.Range("D13").Value = 2050
$ws.Range("E13").Value = 50
$ws.Range("F13").Value = "ELE_NEW_WIND-ON"

# D13 keeps the bottom-row shading used elsewhere for the "2050" cells.
$ws.Range("D10").Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("D13").Value = 2050

# ---------------------------------------------------------------------
# 2. Left-over formatted (but empty) rows below the new bound, mirroring
#    the alternating row shading of the table above (rows 14-15 and
#    20-22), plus a couple of blank spacer rows (12, 19). These are
#    duplicated from the still-unmodified rows 6/7/10 so columns B/D/E/F
#    keep their original alternating shading.
# ---------------------------------------------------------------------

$ws.Rows(12).RowHeight = 13.5

$ws.Range("B6:F6").Copy()
$ws.Range("B14:F14").PasteSpecial(-4122)
$ws.Range("B14:F14").ClearContents()

$ws.Range("B7:F7").Copy()
$ws.Range("B15:F15").PasteSpecial(-4122)
$ws.Range("B15:F15").ClearContents()

$ws.Rows(19).RowHeight = 13.5

$ws.Range("B6:F6").Copy()
$ws.Range("B20:F20").PasteSpecial(-4122)
$ws.Range("B20:F20").ClearContents()

$ws.Range("B7:F7").Copy()
$ws.Range("B21:F21").PasteSpecial(-4122)
$ws.Range("B21:F21").ClearContents()

$ws.Range("B10:F10").Copy()
$ws.Range("B22:F22").PasteSpecial(-4122)
$ws.Range("B22:F22").ClearContents()

# The C column keeps the same (un-shaded) style used by row 13's C cell
# all the way down, regardless of the B/D/E/F alternation.
$ws.Range("C5").Copy()
$ws.Range("C14:C15").PasteSpecial(-4122)
$ws.Range("C20:C22").PasteSpecial(-4122)

$ws.Range("B13:F15").RowHeight = 13.5
$ws.Range("B20:F22").RowHeight = 13.5

# ---------------------------------------------------------------------
# 3. Existing rows 5-10: the bound process changes from the nuclear
#    process (PP_NEW_NUC) to the onshore-wind process (ELE_NEW_WIND-ON),
#    row 5's limit type switches from UP to FX, and the capacity values
#    (col E) are scaled up.
# ---------------------------------------------------------------------

$ws.Range("B5").Value = "FX"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = "ELE_NEW_WIND-ON"

$ws.Range("E6").Value = 4
$ws.Range("F6").Value = "ELE_NEW_WIND-ON"

$ws.Range("E7").Value = 8
$ws.Range("F7").Value = "ELE_NEW_WIND-ON"

$ws.Range("E8").Value = 16
$ws.Range("F8").Value = "ELE_NEW_WIND-ON"

$ws.Range("E9").Value = 16
$ws.Range("F9").Value = "ELE_NEW_WIND-ON"

$ws.Range("E10").Value = 16
$ws.Range("F10").Value = "ELE_NEW_WIND-ON"

# Column F on rows 6-10 picks up the same (unshaded) style already used
# by row 5, so the whole row reads as one visual block.
$ws.Range("F5").Copy()
$ws.Range("F6:F10").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 4. Column F is widened to fit the longer process name, and the
#    selection cursor ends up on G4.
# ---------------------------------------------------------------------

$ws.Columns(6).ColumnWidth = 22.14

$ws.Range("G4").Select()
